# Apply the periodic data refresh ("Update gh-pages to output generated at
# 456a3b4") to all affected sheets: bumped "want to go" counts (column F)
# across 展览/演出/本地生活/全部类型, plus a title annotation on the
# Fun-X event row ("...【免费入场】").
#
# Each update is expressed as an explicit object (Sheet/Cell/Value) rather
# than a nested array, since PowerShell's array-literal flattening would
# otherwise collapse nested @(...) arrays when iterated.

$wb = $excel.ActiveWorkbook

$updates = @(
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F2";  Value = 185 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F4";  Value = 1123 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F5";  Value = 36 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F7";  Value = 14 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "C10"; Value = "杭州·首届Fun-X动漫嘉年华【免费入场】" }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F10"; Value = 317 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F11"; Value = 412 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F12"; Value = 32 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F13"; Value = 303 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F14"; Value = 345 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F15"; Value = 23 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F17"; Value = 404 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F18"; Value = 433 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F19"; Value = 5498 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F21"; Value = 1539 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F22"; Value = 356 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F23"; Value = 4653 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F24"; Value = 115 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F25"; Value = 83 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F26"; Value = 1473 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F29"; Value = 641 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F30"; Value = 30 }
    [PSCustomObject]@{ Sheet = "展览"; Cell = "F32"; Value = 3788 }

    [PSCustomObject]@{ Sheet = "演出"; Cell = "F4";  Value = 105 }

    [PSCustomObject]@{ Sheet = "本地生活"; Cell = "F2"; Value = 9370 }
    [PSCustomObject]@{ Sheet = "本地生活"; Cell = "F4"; Value = 2116 }

    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F2";  Value = 9370 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F4";  Value = 2116 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F5";  Value = 185 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F7";  Value = 1123 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F8";  Value = 36 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F10"; Value = 14 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "C13"; Value = "杭州·首届Fun-X动漫嘉年华【免费入场】" }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F13"; Value = 317 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F14"; Value = 412 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F15"; Value = 32 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F16"; Value = 303 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F17"; Value = 345 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F18"; Value = 23 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F23"; Value = 404 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F24"; Value = 433 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F25"; Value = 5499 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F27"; Value = 1539 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F30"; Value = 356 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F32"; Value = 4653 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F33"; Value = 115 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F34"; Value = 83 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F35"; Value = 1473 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F38"; Value = 641 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F39"; Value = 30 }
    [PSCustomObject]@{ Sheet = "全部类型"; Cell = "F46"; Value = 3788 }
)

foreach ($item in $updates) {
    $ws = $wb.Worksheets.Item($item.Sheet)
    $ws.Range($item.Cell).Value = $item.Value
}
